$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2090.818
$ws.Range("I19").Value = 4590.5
$ws.Range("J19").Value = 662.4286
$ws.Range("K19").Value = 4590.5
$ws.Range("L19").Value = 662.4286
$ws.Range("M19").Value = -4415.5
$ws.Range("N19").Value = -1012.4286

$ws.Range("H116").Value = 2555.9375
$ws.Range("I116").Value = 2641.3
$ws.Range("J116").Value = 2413.6667
$ws.Range("K116").Value = 2641.3
$ws.Range("L116").Value = 2413.6667
$ws.Range("M116").Value = 800.6999999999998
$ws.Range("N116").Value = -9297.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1425
$ws.Range("I45").Value = 1319.4286
$ws.Range("J45").Value = 1671.3334
$ws.Range("K45").Value = 1319.4286
$ws.Range("L45").Value = 1671.3334
$ws.Range("M45").Value = -942.4286
$ws.Range("N45").Value = -2425.3334

$ws.Range("H63").Value = 6044
$ws.Range("I63").Value = 6426.2
$ws.Range("J63").Value = 2222
$ws.Range("K63").Value = 6426.2
$ws.Range("L63").Value = 2222
$ws.Range("M63").Value = -5740.2
$ws.Range("N63").Value = -3594

$ws.Range("H66").Value = 6044
$ws.Range("I66").Value = 6426.2
$ws.Range("J66").Value = 2222
$ws.Range("K66").Value = 32131
$ws.Range("L66").Value = 11110
$ws.Range("M66").Value = -28699
$ws.Range("N66").Value = -17974

$ws.Range("H74").Value = 197872.55
$ws.Range("I74").Value = 233312.77
$ws.Range("J74").Value = 80647.234
$ws.Range("K74").Value = 233312.77
$ws.Range("L74").Value = 80647.234
$ws.Range("M74").Value = -232438.77
$ws.Range("N74").Value = -82395.234

$ws.Range("H77").Value = 197872.55
$ws.Range("I77").Value = 233312.77
$ws.Range("J77").Value = 80647.234
$ws.Range("K77").Value = 1166563.85
$ws.Range("L77").Value = 403236.17
$ws.Range("M77").Value = -1162195.85
$ws.Range("N77").Value = -411972.17

$ws.Range("H102").Value = 4454.5454
$ws.Range("I102").Value = 1666.6666
$ws.Range("K102").Value = 1666.6666
$ws.Range("M102").Value = -44.66660000000002

$ws.Range("H110").Value = 1131.0667
$ws.Range("I110").Value = 1039.9714
$ws.Range("J110").Value = 1449.9
$ws.Range("K110").Value = 1039.9714
$ws.Range("L110").Value = 1449.9
$ws.Range("M110").Value = 1005.0286
$ws.Range("N110").Value = -5539.9

$ws.Range("H132").Value = 31521.223
$ws.Range("I132").Value = 45427.75
$ws.Range("K132").Value = 136283.25
$ws.Range("M132").Value = -133753.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 833.0769
$ws.Range("I107").Value = 719.3939
$ws.Range("J107").Value = 1458.3334
$ws.Range("K107").Value = 719.3939
$ws.Range("L107").Value = 1458.3334
$ws.Range("M107").Value = 1200.6061
$ws.Range("N107").Value = -5298.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1859.4386
$ws.Range("I31").Value = 870.0909
$ws.Range("J31").Value = 5208
$ws.Range("K31").Value = 870.0909
$ws.Range("L31").Value = 5208
$ws.Range("M31").Value = -575.0909
$ws.Range("N31").Value = -5798

$ws.Range("H34").Value = 1859.4386
$ws.Range("I34").Value = 870.0909
$ws.Range("J34").Value = 5208
$ws.Range("K34").Value = 870.0909
$ws.Range("L34").Value = 5208
$ws.Range("M34").Value = -668.0909
$ws.Range("N34").Value = -5612

$ws.Range("H58").Value = 8757
$ws.Range("I58").Value = 50756
$ws.Range("J58").Value = 2757.1428
$ws.Range("K58").Value = 50756
$ws.Range("L58").Value = 2757.1428
$ws.Range("M58").Value = -50553
$ws.Range("N58").Value = -3163.1428

$ws.Range("H107").Value = 318.36365
$ws.Range("I107").Value = 334.63333
$ws.Range("J107").Value = 283.5
$ws.Range("K107").Value = 334.63333
$ws.Range("L107").Value = 283.5
$ws.Range("M107").Value = 1585.36667
$ws.Range("N107").Value = -4123.5

$ws.Range("H132").Value = 2485.3928
$ws.Range("I132").Value = 1644.5714
$ws.Range("J132").Value = 3326.2144
$ws.Range("K132").Value = 4933.7142
$ws.Range("L132").Value = 9978.643199999999
$ws.Range("M132").Value = -2403.7142
$ws.Range("N132").Value = -15038.6432

$ws.Range("H136").Value = 8757
$ws.Range("I136").Value = 50756
$ws.Range("J136").Value = 2757.1428
$ws.Range("K136").Value = 152268
$ws.Range("L136").Value = 8271.428400000001
$ws.Range("M136").Value = -149718
$ws.Range("N136").Value = -13371.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1012.7955
$ws.Range("I5").Value = 409.92
$ws.Range("K5").Value = 1229.76
$ws.Range("M5").Value = -1117.76

$ws.Range("H12").Value = 40.322582
$ws.Range("J12").Value = 45.913044
$ws.Range("L12").Value = 137.739132
$ws.Range("N12").Value = -483.739132

$ws.Range("H113").Value = 633.125
$ws.Range("I113").Value = 619.45
$ws.Range("J113").Value = 701.5
$ws.Range("K113").Value = 1858.35
$ws.Range("L113").Value = 2104.5
$ws.Range("M113").Value = 311.6499999999999
$ws.Range("N113").Value = -6444.5

$ws.Range("H116").Value = 4756
$ws.Range("I116").Value = 300
$ws.Range("J116").Value = 5784.3076
$ws.Range("K116").Value = 900
$ws.Range("L116").Value = 17352.9228
$ws.Range("M116").Value = 2542
$ws.Range("N116").Value = -24236.9228

$ws.Range("H122").Value = 564.24
$ws.Range("I122").Value = 406.33334
$ws.Range("J122").Value = 970.2857
$ws.Range("K122").Value = 3657.00006
$ws.Range("L122").Value = 8732.5713
$ws.Range("M122").Value = -1207.00006
$ws.Range("N122").Value = -13632.5713

$ws.Range("H135").Value = 1012.7955
$ws.Range("I135").Value = 409.92
$ws.Range("K135").Value = 3689.28
$ws.Range("M135").Value = -1154.28

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4510.5835
$ws.Range("I80").Value = 5544.3184
$ws.Range("J80").Value = 2886.1428
$ws.Range("K80").Value = 5544.3184
$ws.Range("L80").Value = 2886.1428
$ws.Range("M80").Value = -4546.3184
$ws.Range("N80").Value = -4882.1428

$ws.Range("H83").Value = 4510.5835
$ws.Range("I83").Value = 5544.3184
$ws.Range("J83").Value = 2886.1428
$ws.Range("K83").Value = 27721.592
$ws.Range("L83").Value = 14430.714
$ws.Range("M83").Value = -22729.592
$ws.Range("N83").Value = -24414.714

$ws.Range("H126").Value = 1974.3889
$ws.Range("I126").Value = 1758.0741
$ws.Range("K126").Value = 5274.2223
$ws.Range("M126").Value = -2804.2223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2291.6667
$ws.Range("I122").Value = 1875
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5625
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3175
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1786.8158
$ws.Range("I132").Value = 944.4815
$ws.Range("J132").Value = 3854.3635
$ws.Range("K132").Value = 2833.4445
$ws.Range("L132").Value = 11563.0905
$ws.Range("M132").Value = -303.4445000000001
$ws.Range("N132").Value = -16623.0905

$ws.Range("H136").Value = 14043547
$ws.Range("I136").Value = 23280544
$ws.Range("J136").Value = 347310
$ws.Range("K136").Value = 69841632
$ws.Range("L136").Value = 1041930
$ws.Range("M136").Value = -69839082
$ws.Range("N136").Value = -1047030
